$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC row 11 ---
$ws_ALC.Range("H11").Value2 = 58109.383
$ws_ALC.Range("I11").Value2 = 58109.383
$ws_ALC.Range("K11").Value2 = 58109.383
$ws_ALC.Range("M11").Value2 = -57969.383

# --- ALC row 19 ---
$ws_ALC.Range("H19").Value2 = 1061.5
$ws_ALC.Range("I19").Value2 = 1697.75
$ws_ALC.Range("J19").Value2 = 425.25
$ws_ALC.Range("K19").Value2 = 1697.75
$ws_ALC.Range("L19").Value2 = 425.25
$ws_ALC.Range("M19").Value2 = -1522.75
$ws_ALC.Range("N19").Value2 = -775.25

# --- ALC row 49 ---
$ws_ALC.Range("H49").Value2 = 301
$ws_ALC.Range("I49").Value2 = 307
$ws_ALC.Range("J49").Value2 = 299.5
$ws_ALC.Range("K49").Value2 = 921
$ws_ALC.Range("L49").Value2 = 898.5
$ws_ALC.Range("M49").Value2 = -785
$ws_ALC.Range("N49").Value2 = -1170.5

# --- ALC row 58 ---
$ws_ALC.Range("H58").Value2 = 1166.5555
$ws_ALC.Range("I58").Value2 = 54.142857
$ws_ALC.Range("K58").Value2 = 162.428571
$ws_ALC.Range("M58").Value2 = -12.42857100000001

# --- ALC row 80 ---
$ws_ALC.Range("H80").Value2 = 3856.2856
$ws_ALC.Range("I80").Value2 = 2598.8
$ws_ALC.Range("K80").Value2 = 7796.400000000001
$ws_ALC.Range("M80").Value2 = -6798.400000000001

# --- ALC row 83 ---
$ws_ALC.Range("H83").Value2 = 3856.2856
$ws_ALC.Range("I83").Value2 = 2598.8
$ws_ALC.Range("K83").Value2 = 23389.2
$ws_ALC.Range("M83").Value2 = -18397.2

# --- ALC row 116 ---
$ws_ALC.Range("H116").Value2 = 24267.285
$ws_ALC.Range("I116").Value2 = 6119
$ws_ALC.Range("K116").Value2 = 6119
$ws_ALC.Range("M116").Value2 = -2677

# --- ALC row 138 ---
$ws_ALC.Range("H138").Value2 = 2955.1226
$ws_ALC.Range("I138").Value2 = 1925.2354
$ws_ALC.Range("J138").Value2 = 3502.25
$ws_ALC.Range("K138").Value2 = 5775.706200000001
$ws_ALC.Range("L138").Value2 = 10506.75
$ws_ALC.Range("M138").Value2 = -635.7062000000005
$ws_ALC.Range("N138").Value2 = -20786.75

# --- ARM row 61 ---
$ws_ARM.Range("H61").Value2 = 1473570.2
$ws_ARM.Range("I61").Value2 = 2820.1833
$ws_ARM.Range("K61").Value2 = 2820.1833
$ws_ARM.Range("M61").Value2 = -2608.1833

# --- ARM row 74 ---
$ws_ARM.Range("H74").Value2 = 1055430.1
$ws_ARM.Range("I74").Value2 = 1392261.9
$ws_ARM.Range("J74").Value2 = 19024.691
$ws_ARM.Range("K74").Value2 = 1392261.9
$ws_ARM.Range("L74").Value2 = 19024.691
$ws_ARM.Range("M74").Value2 = -1391387.9
$ws_ARM.Range("N74").Value2 = -20772.691

# --- ARM row 77 ---
$ws_ARM.Range("H77").Value2 = 1055430.1
$ws_ARM.Range("I77").Value2 = 1392261.9
$ws_ARM.Range("J77").Value2 = 19024.691
$ws_ARM.Range("K77").Value2 = 6961309.5
$ws_ARM.Range("L77").Value2 = 95123.45499999999
$ws_ARM.Range("M77").Value2 = -6956941.5
$ws_ARM.Range("N77").Value2 = -103859.455

# --- ARM row 122 ---
$ws_ARM.Range("H122").Value2 = 1984.7667
$ws_ARM.Range("I122").Value2 = 1924.9231
$ws_ARM.Range("J122").Value2 = 2373.75
$ws_ARM.Range("K122").Value2 = 5774.7693
$ws_ARM.Range("L122").Value2 = 7121.25
$ws_ARM.Range("M122").Value2 = -3324.7693
$ws_ARM.Range("N122").Value2 = -12021.25

# --- ARM row 132 ---
$ws_ARM.Range("H132").Value2 = 584892.4
$ws_ARM.Range("I132").Value2 = 644440.4399999999
$ws_ARM.Range("K132").Value2 = 1933321.32
$ws_ARM.Range("M132").Value2 = -1930791.32

# --- ARM row 136 ---
$ws_ARM.Range("H136").Value2 = 1473570.2
$ws_ARM.Range("I136").Value2 = 2820.1833
$ws_ARM.Range("K136").Value2 = 8460.5499
$ws_ARM.Range("M136").Value2 = -5910.5499

# --- BSM row 124 ---
$ws_BSM.Range("H124").Value2 = 50000
$ws_BSM.Range("J124").Value2 = 50000
$ws_BSM.Range("L124").Value2 = 50000
$ws_BSM.Range("N124").Value2 = -59820

# --- CRP row 31 ---
$ws_CRP.Range("H31").Value2 = 2227086.5
$ws_CRP.Range("I31").Value2 = 3475529.2
$ws_CRP.Range("K31").Value2 = 3475529.2
$ws_CRP.Range("M31").Value2 = -3475234.2

# --- CRP row 34 ---
$ws_CRP.Range("H34").Value2 = 2227086.5
$ws_CRP.Range("I34").Value2 = 3475529.2
$ws_CRP.Range("K34").Value2 = 3475529.2
$ws_CRP.Range("M34").Value2 = -3475327.2

# --- CRP row 105 ---
$ws_CRP.Range("H105").Value2 = 6569.5835
$ws_CRP.Range("I105").Value2 = 7231.0625
$ws_CRP.Range("K105").Value2 = 7231.0625
$ws_CRP.Range("M105").Value2 = -5484.0625

# --- CRP row 132 ---
$ws_CRP.Range("H132").Value2 = 1849.7435
$ws_CRP.Range("I132").Value2 = 1501.4062
$ws_CRP.Range("K132").Value2 = 4504.2186
$ws_CRP.Range("M132").Value2 = -1974.2186

# --- CRP row 134 ---
$ws_CRP.Range("H134").Value2 = 2074.0908
$ws_CRP.Range("I134").Value2 = 1366.258
$ws_CRP.Range("K134").Value2 = 4098.774
$ws_CRP.Range("M134").Value2 = -1563.774

# --- CUL row 12 ---
$ws_CUL.Range("H12").Value2 = 310.64706
$ws_CUL.Range("I12").Value2 = 636
$ws_CUL.Range("J12").Value2 = 21.444445
$ws_CUL.Range("K12").Value2 = 1908
$ws_CUL.Range("L12").Value2 = 64.33333500000001
$ws_CUL.Range("M12").Value2 = -1735
$ws_CUL.Range("N12").Value2 = -410.333335

# --- CUL row 34 ---
$ws_CUL.Range("H34").Value2 = 5000
$ws_CUL.Range("J34").Value2 = 5000
$ws_CUL.Range("L34").Value2 = 15000
$ws_CUL.Range("N34").Value2 = -15168

# --- CUL row 39 ---
$ws_CUL.Range("H39").Value2 = 3457
$ws_CUL.Range("J39").Value2 = 4599.8
$ws_CUL.Range("L39").Value2 = 13799.4
$ws_CUL.Range("N39").Value2 = -14387.4

# --- CUL row 55 ---
$ws_CUL.Range("H55").Value2 = 76366100
$ws_CUL.Range("I55").Value2 = 93334670
$ws_CUL.Range("K55").Value2 = 280004010
$ws_CUL.Range("M55").Value2 = -280003833

# --- CUL row 81 ---
$ws_CUL.Range("H81").Value2 = 1666.6666

# --- CUL row 84 ---
$ws_CUL.Range("H84").Value2 = 1666.6666

# --- CUL row 109 ---
$ws_CUL.Range("H109").Value2 = 6662.3687
$ws_CUL.Range("I109").Value2 = 4070.5557
$ws_CUL.Range("J109").Value2 = 8995
$ws_CUL.Range("K109").Value2 = 12211.6671
$ws_CUL.Range("L109").Value2 = 26985
$ws_CUL.Range("M109").Value2 = -11171.6671
$ws_CUL.Range("N109").Value2 = -29065

# --- CUL row 132 ---
$ws_CUL.Range("H132").Value2 = 8857.357
$ws_CUL.Range("I132").Value2 = 1338
$ws_CUL.Range("J132").Value2 = 13034.777
$ws_CUL.Range("K132").Value2 = 12042
$ws_CUL.Range("L132").Value2 = 117312.993
$ws_CUL.Range("M132").Value2 = -9512
$ws_CUL.Range("N132").Value2 = -122372.993

# --- CUL row 138 ---
$ws_CUL.Range("H138").Value2 = 103210
$ws_CUL.Range("J138").Value2 = 0
$ws_CUL.Range("L138").Value2 = 0
$ws_CUL.Range("N138").ClearContents()

# --- GSM row 122 ---
$ws_GSM.Range("H122").Value2 = 64237.06
$ws_GSM.Range("I122").Value2 = 76865.86
$ws_GSM.Range("K122").Value2 = 230597.58
$ws_GSM.Range("M122").Value2 = -228147.58

# --- GSM row 123 ---
$ws_GSM.Range("H123").Value2 = 45898
$ws_GSM.Range("J123").Value2 = 45898
$ws_GSM.Range("L123").Value2 = 45898
$ws_GSM.Range("N123").Value2 = -50798

# --- GSM row 132 ---
$ws_GSM.Range("H132").Value2 = 11667.521
$ws_GSM.Range("I132").Value2 = 10259.395
$ws_GSM.Range("K132").Value2 = 30778.185
$ws_GSM.Range("M132").Value2 = -28248.185

# --- LTW row 7 ---
$ws_LTW.Range("H7").Value2 = 5776.8125
$ws_LTW.Range("I7").Value2 = 6000
$ws_LTW.Range("K7").Value2 = 6000
$ws_LTW.Range("M7").Value2 = -5888

# --- LTW row 40 ---
$ws_LTW.Range("H40").Value2 = 6600.5186
$ws_LTW.Range("I40").Value2 = 6328.6
$ws_LTW.Range("J40").Value2 = 9999.5
$ws_LTW.Range("K40").Value2 = 6328.6
$ws_LTW.Range("L40").Value2 = 9999.5
$ws_LTW.Range("M40").Value2 = -6192.6
$ws_LTW.Range("N40").Value2 = -10271.5

# --- LTW row 69 ---
$ws_LTW.Range("H69").Value2 = 9499.333000000001
$ws_LTW.Range("J69").Value2 = 9499.333000000001
$ws_LTW.Range("L69").Value2 = 9499.333000000001
$ws_LTW.Range("N69").Value2 = -11121.333

# --- LTW row 72 ---
$ws_LTW.Range("H72").Value2 = 9499.333000000001
$ws_LTW.Range("J72").Value2 = 9499.333000000001
$ws_LTW.Range("L72").Value2 = 28497.999
$ws_LTW.Range("N72").Value2 = -36609.999

# --- LTW row 82 ---
$ws_LTW.Range("H82").Value2 = 1615.68
$ws_LTW.Range("I82").Value2 = 1458.7727
$ws_LTW.Range("J82").Value2 = 2766.3333
$ws_LTW.Range("K82").Value2 = 1458.7727
$ws_LTW.Range("L82").Value2 = 2766.3333
$ws_LTW.Range("M82").Value2 = -1097.7727
$ws_LTW.Range("N82").Value2 = -3488.3333

# --- LTW row 85 ---
$ws_LTW.Range("H85").Value2 = 1615.68
$ws_LTW.Range("I85").Value2 = 1458.7727
$ws_LTW.Range("J85").Value2 = 2766.3333
$ws_LTW.Range("K85").Value2 = 1458.7727
$ws_LTW.Range("L85").Value2 = 2766.3333
$ws_LTW.Range("M85").Value2 = -210.7727
$ws_LTW.Range("N85").Value2 = -5262.3333

# --- LTW row 126 ---
$ws_LTW.Range("H126").Value2 = 5776.8125
$ws_LTW.Range("I126").Value2 = 6000
$ws_LTW.Range("K126").Value2 = 18000
$ws_LTW.Range("M126").Value2 = -15530

# --- LTW row 132 ---
$ws_LTW.Range("H132").Value2 = 3654624.5
$ws_LTW.Range("I132").Value2 = 6494276
$ws_LTW.Range("J132").Value2 = 3644.0715
$ws_LTW.Range("K132").Value2 = 19482828
$ws_LTW.Range("L132").Value2 = 10932.2145
$ws_LTW.Range("M132").Value2 = -19480298
$ws_LTW.Range("N132").Value2 = -15992.2145

# --- WVR row 105 ---
$ws_WVR.Range("H105").Value2 = 60000
$ws_WVR.Range("J105").Value2 = 60000
$ws_WVR.Range("L105").Value2 = 60000
$ws_WVR.Range("N105").Value2 = -66988

# --- WVR row 122 ---
$ws_WVR.Range("H122").Value2 = 41438.965
$ws_WVR.Range("I122").Value2 = 957.2105
$ws_WVR.Range("K122").Value2 = 2871.6315
$ws_WVR.Range("M122").Value2 = -421.6315

# --- WVR row 132 ---
$ws_WVR.Range("H132").Value2 = 2925294
$ws_WVR.Range("I132").Value2 = 3334511.2
$ws_WVR.Range("K132").Value2 = 10003533.6
$ws_WVR.Range("M132").Value2 = -10001003.6
